# Fruta / hortaliza, semanal
# Insert two new weekly records at rows 58-59 (pushing the existing rows
# 58..71 down to 60..73), then populate the two new rows with the new
# observations.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 58:59 - everything currently at row 58 and below
# shifts down by two rows (old row 58 -> new row 60, ..., old row 71 -> new row 73).
$ws.Range("A58:A59").EntireRow.Insert()

# New row 58: Artic Pride, Primera, Region de Coquimbo
$ws.Cells.Item(58, 1).Value = 1
$ws.Cells.Item(58, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(58, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(58, 4).Value = 44917
$ws.Cells.Item(58, 5).Value = 15
$ws.Cells.Item(58, 6).Value = "Fruta"
$ws.Cells.Item(58, 7).Value = 100103
$ws.Cells.Item(58, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(58, 9).Value = 100103006
$ws.Cells.Item(58, 10).Value = "Nectarín"
$ws.Cells.Item(58, 11).Value = "Artic Pride"
$ws.Cells.Item(58, 12).Value = "Primera"
$ws.Cells.Item(58, 13).Value = 300
$ws.Cells.Item(58, 14).Value = 20000
$ws.Cells.Item(58, 15).Value = 21000
$ws.Cells.Item(58, 16).Value = 20500
$ws.Cells.Item(58, 17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item(58, 18).Value = "Región de Coquimbo"
$ws.Cells.Item(58, 19).Value = 1139
$ws.Cells.Item(58, 20).Value = 18

# New row 59: Super Queen, Primera, Region de Coquimbo
$ws.Cells.Item(59, 1).Value = 1
$ws.Cells.Item(59, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(59, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(59, 4).Value = 44917
$ws.Cells.Item(59, 5).Value = 15
$ws.Cells.Item(59, 6).Value = "Fruta"
$ws.Cells.Item(59, 7).Value = 100103
$ws.Cells.Item(59, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(59, 9).Value = 100103006
$ws.Cells.Item(59, 10).Value = "Nectarín"
$ws.Cells.Item(59, 11).Value = "Super Queen"
$ws.Cells.Item(59, 12).Value = "Primera"
$ws.Cells.Item(59, 13).Value = 400
$ws.Cells.Item(59, 14).Value = 20000
$ws.Cells.Item(59, 15).Value = 21000
$ws.Cells.Item(59, 16).Value = 20625
$ws.Cells.Item(59, 17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item(59, 18).Value = "Región de Coquimbo"
$ws.Cells.Item(59, 19).Value = 1146
$ws.Cells.Item(59, 20).Value = 18

# Make sure the date cells use the same date/time number format as the
# rest of the "Fecha" column.
$ws.Range("D58:D59").NumberFormat = "YYYY-MM-DD HH:MM:SS"
